$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J5").Value = "increased attack speed (lower timebetweenattacks) + lifesteal + red shader"
$ws.Range("I5").Value = "lasts 10s"

$ws.Range("I5").Select()
